$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (AMM)
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 1222
$ws.Range("D2").Value = 1499
$ws.Range("E2").Value = 106
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 21
$ws.Range("H2").Value = 159
$ws.Range("I2").Value = 1595.1
$ws.Range("J2").Value = -6.024700645727532

# Row 3 (IPR)
$ws.Range("C3").Value = 163
$ws.Range("D3").Value = 165
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = 158
$ws.Range("J3").Value = 4.430379746835444

# Row 4 (MIG)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 342
$ws.Range("D4").Value = 353
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 9
$ws.Range("J4").Value = 253

# Row 5 (MOB)
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 875
$ws.Range("D5").Value = 1012
$ws.Range("E5").Value = 61
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 64
$ws.Range("I5").Value = 1028
$ws.Range("J5").Value = -1.556420233463029

# Row 6 (MOB PRE)
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 711
$ws.Range("D6").Value = 854
$ws.Range("E6").Value = 136
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 4
$ws.Range("I6").Value = 758
$ws.Range("J6").Value = 12.66490765171504

# Row 7 (MSK)
$ws.Range("C7").Value = 212
$ws.Range("D7").Value = 237
$ws.Range("E7").Value = 24
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("I7").Value = 232
$ws.Range("J7").Value = 2.155172413793105

# Row 9 (TEC)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 302
$ws.Range("D9").Value = 305
$ws.Range("E9").Value = 2
$ws.Range("G9").Value = 4
$ws.Range("I9").Value = 821
$ws.Range("J9").Value = -62.85018270401949

# Row 10 (TST)
$ws.Range("C10").Value = 73
$ws.Range("D10").Value = 80
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 1
$ws.Range("I10").Value = 107
$ws.Range("J10").Value = -25.23364485981309

# Row 11 (VIP)
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 2
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = -40

# Row 12 (WLC)
$ws.Range("C12").Value = 36
$ws.Range("D12").Value = 37
$ws.Range("E12").Value = 1
$ws.Range("I12").Value = 74
$ws.Range("J12").Value = -50
